$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.320.15"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "1.953.14"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'243.45"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").Value = "'58.40"
$ws.Range("E7").Value = "  -7.89%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.371"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").Value = "'55.68"
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("D11").Value = "'0.0837"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'21.85"
$ws.Range("E13").Value = "  -7.11%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.831"
$ws.Range("E14").Value = "  -8.66%  "
$ws.Range("D15").Value = "2.242.26"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "'13.51"
$ws.Range("E16").Value = "  -6.57%  "
$ws.Range("D17").Value = "'5.33"
$ws.Range("E17").Value = "  -4.03%  "
$ws.Range("D18").Value = "1.996.34"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "36.276.18"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'70.07"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0875"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'229.99"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  -6.95%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "'9.50"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "'163.92"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").Value = "'19.64"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  -10.89%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "'4.72"
$ws.Range("E33").Value = "  -7.14%  "
$ws.Range("D34").Value = "'0.0636"
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("D35").Value = "'4.31"
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.11"
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  -9.50%  "
$ws.Range("D40").Value = "'2.89"
$ws.Range("E40").Value = "  -7.50%  "
$ws.Range("D41").Value = "'0.0979"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.87"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.18"
$ws.Range("E43").Value = "  -6.21%  "
$ws.Range("D44").Value = "'0.0210"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.04"
$ws.Range("E45").Value = "  -8.41%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'15.69"
$ws.Range("E46").Value = "  -8.72%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'88.74"
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.37"
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "1.344.38"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "'45.10"
$ws.Range("E51").Value = "  -1.18%  "
